$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 67

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/06"
$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 5
$ws.Cells.Item($row, 4).Value = 71
